$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 166, shifting existing rows 166.. down to 167..
$ws.Rows("166").Insert()

# Populate the newly inserted row 166 with the new record
$ws.Range("A166").Value = 6
$ws.Range("B166").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C166").Value = "Metropolitana"
$ws.Range("D166").Value = 44777
$ws.Range("E166").Value = 13
$ws.Range("F166").Value = 100112026
$ws.Range("G166").Value = "Haba"
$ws.Range("H166").Value = "Sin especificar"
$ws.Range("I166").Value = "Primera"
$ws.Range("J166").Value = 400
$ws.Range("K166").Value = 17000
$ws.Range("L166").Value = 18000
$ws.Range("M166").Value = 17425
$ws.Range("N166").Value = "$/saco 25 kilos"
$ws.Range("O166").Value = "Región de Coquimbo"
$ws.Range("P166").Value = 697
$ws.Range("Q166").Value = 25
$ws.Range("R166").Value = "Hortaliza"
